$wb = $excel.ActiveWorkbook

# --- Update the "SoCDTtiNTY-frgt" sheet (freight) ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Row 2 ("LDVs"): B2 gets its own formula, C2:H2 become a shared-formula block.
$wsFrgt.Range("B2").Formula = "=1/17"
$wsFrgt.Range("C2:H2").Formula = "=1/17"

# Row 3 ("HDVs"): B3 gets its own formula, C3:H3 become a shared-formula block.
$wsFrgt.Range("B3").Formula = "=1/19"
$wsFrgt.Range("C3:H3").Formula = "=1/19"

# Update the visible selection on that sheet to span B2:H3, anchored at B2.
$wsFrgt.Activate() | Out-Null
$wsFrgt.Range("B2:H3").Select() | Out-Null

# --- Make "About" the active/selected sheet (tab) in the workbook ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
